$wb = $excel.ActiveWorkbook

function Add-LogRows($tgt, $rws) {
  foreach ($row in $rws) {
    $r = $row[0]
    for ($c = 1; $c -le 6; $c++) {
      $cell = $tgt.Cells.Item($r, $c)
      $cell.NumberFormat = "@"
      $cell.Value = $row[$c]
    }
  }
}

$ws = $wb.Worksheets.Item("PIR")
$PIRRows = @(
  @("31","2026-02-04","14:03:13","14:00","Bathroom","No Motion","Inactive"),
  @("32","2026-02-04","14:03:18","14:00","Bathroom","No Motion","Inactive"),
  @("33","2026-02-04","14:03:23","14:00","Bathroom","Motion Detected","Active"),
  @("34","2026-02-04","14:03:30","14:00","Bathroom","No Motion","Inactive"),
  @("35","2026-02-04","14:03:35","14:00","Bathroom","No Motion","Inactive"),
  @("36","2026-02-04","14:03:40","14:00","Bathroom","No Motion","Inactive"),
  @("37","2026-02-04","14:03:45","14:00","Bathroom","No Motion","Inactive"),
  @("38","2026-02-04","14:03:49","14:00","Bathroom","Motion Detected","Active"),
  @("39","2026-02-04","14:03:57","14:00","Bathroom","No Motion","Inactive"),
  @("40","2026-02-04","14:04:02","14:00","Bathroom","No Motion","Inactive"),
  @("41","2026-02-04","14:04:04","14:00","Bathroom","Motion Detected","Active")
)
Add-LogRows $ws $PIRRows

$ws = $wb.Worksheets.Item("Humidity")
$HumidityRows = @(
  @("25","2026-02-04","14:03:13","14:00","Bathroom","77.5%","Active"),
  @("26","2026-02-04","14:03:17","14:00","Bathroom","78.4%","Active"),
  @("27","2026-02-04","14:03:22","14:00","Bathroom","77.2%","Active"),
  @("28","2026-02-04","14:03:27","14:00","Bathroom","77.9%","Active"),
  @("29","2026-02-04","14:03:38","14:00","Bathroom","77.8%","Active"),
  @("30","2026-02-04","14:03:43","14:00","Bathroom","76.7%","Active"),
  @("31","2026-02-04","14:03:48","14:00","Bathroom","77.5%","Active"),
  @("32","2026-02-04","14:03:53","14:00","Bathroom","76.6%","Active"),
  @("33","2026-02-04","14:03:58","14:00","Bathroom","77.6%","Active"),
  @("34","2026-02-04","14:04:03","14:00","Bathroom","76.7%","Active")
)
Add-LogRows $ws $HumidityRows

$ws = $wb.Worksheets.Item("Temperature")
$TemperatureRows = @(
  @("25","2026-02-04","14:03:13","14:00","Bathroom","24.8C","Active"),
  @("26","2026-02-04","14:03:18","14:00","Bathroom","24.8C","Active"),
  @("27","2026-02-04","14:03:23","14:00","Bathroom","24.8C","Active"),
  @("28","2026-02-04","14:03:28","14:00","Bathroom","24.8C","Active"),
  @("29","2026-02-04","14:03:38","14:00","Bathroom","24.8C","Active"),
  @("30","2026-02-04","14:03:43","14:00","Bathroom","24.8C","Active"),
  @("31","2026-02-04","14:03:48","14:00","Bathroom","24.8C","Active"),
  @("32","2026-02-04","14:03:53","14:00","Bathroom","24.8C","Active"),
  @("33","2026-02-04","14:03:58","14:00","Bathroom","24.8C","Active"),
  @("34","2026-02-04","14:04:03","14:00","Bathroom","24.8C","Active")
)
Add-LogRows $ws $TemperatureRows
